$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 368, pushing the existing
# rows 368-393 down to 369-394 (dimension grows from A1:T393 to A1:T394).
$ws.Rows.Item(368).Insert()

$ws.Cells.Item(368, 1).Value = 6
$ws.Cells.Item(368, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(368, 3).Value = "Metropolitana"
$ws.Cells.Item(368, 4).Value = 44826
$ws.Cells.Item(368, 5).Value = 13
$ws.Cells.Item(368, 6).Value = "Fruta"
$ws.Cells.Item(368, 7).Value = 100101
$ws.Cells.Item(368, 8).Value = "Berries"
$ws.Cells.Item(368, 9).Value = 100101001
$ws.Cells.Item(368, 10).Value = "Arándano (blue)"
$ws.Cells.Item(368, 11).Value = "Sin especificar"
$ws.Cells.Item(368, 12).Value = "Primera"
$ws.Cells.Item(368, 13).Value = 350
$ws.Cells.Item(368, 14).Value = 6000
$ws.Cells.Item(368, 15).Value = 7000
$ws.Cells.Item(368, 16).Value = 6500
$ws.Cells.Item(368, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(368, 18).Value = "Perú"
$ws.Cells.Item(368, 19).Value = 4333
$ws.Cells.Item(368, 20).Value = 1.5
